# Add a new "ShopForPaintings" worksheet at the end of the workbook,
# mirroring the structure of the existing "ShopFor..." sheets (e.g.
# "ShopForChargers"), then make it the active/selected sheet.

$wb = $excel.ActiveWorkbook

$srcSheet = $wb.Worksheets.Item("ShopForChargers")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "ShopForPaintings"

# Match the default row height (16) used by the sibling "ShopFor..." sheets
# (which comes from the workbook's size-12 default font) instead of this
# new sheet's size-11 default.
$newSheet.StandardHeight = 16

# Copy the cell formatting (font/number-format) from an existing
# "id" cell (A2, style index 3) and an existing "description" cell
# (D2, style index 2) so the new cells pick up the same styles
# instead of minting new ones.
$srcSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)   # xlPasteFormats
$srcSheet.Range("D2").Copy()
$newSheet.Range("B2").PasteSpecial(-4122)   # xlPasteFormats

# Fill in the data row.
$newSheet.Range("A2").Value = "8"
$newSheet.Range("B2").Value = "Watercolor Art Paintings"

# Widen column B to fit the description text (matches the other
# "ShopFor" sheets' pattern of widening the last data column).
$newSheet.Columns.Item(2).ColumnWidth = 28.25

# Select B2 (the last populated cell) as the active cell on the new sheet.
$newSheet.Range("B2").Select()

# Make the new sheet the active tab.
$newSheet.Activate()
